$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '60.306.41'
Set-TextValue 'E2' '  -5.92%  '
Set-TextValue 'D3' '3.297.81'
Set-TextValue 'E3' '  -4.98%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '560.83'
Set-TextValue 'E5' '  -3.89%  '
Set-TextValue 'D6' '128.14'
Set-TextValue 'E6' '  -2.25%  '
Set-TextValue 'E7' '  -0.02%  '
Set-TextValue 'D8' '3.293.01'
Set-TextValue 'E8' '  -5.09%  '
Set-TextValue 'E9' '  -1.75%  '
Set-TextValue 'E10' '  -3.90%  '
Set-TextValue 'E11' '  -4.71%  '
Set-TextValue 'D12' '0.371'
Set-TextValue 'E12' '  -3.98%  '
Set-TextValue 'D13' '3.851.47'
Set-TextValue 'E13' '  -5.28%  '
Set-TextValue 'D15' '3.283.96'
Set-TextValue 'E15' '  -5.44%  '
Set-TextValue 'E16' '  -5.61%  '
Set-TextValue 'B17' 'Avalanche'
Set-TextValue 'C17' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D17' '24.17'
Set-TextValue 'E17' '  -0.32%  '
Set-TextValue 'B18' 'WrappedBTC'
Set-TextValue 'C18' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D18' '60.463.11'
Set-TextValue 'E18' '  -5.60%  '
Set-TextValue 'D19' '5.64'
Set-TextValue 'E19' '  -0.51%  '
Set-TextValue 'D20' '13.31'
Set-TextValue 'E20' '  -0.52%  '
Set-TextValue 'D21' '9.06'
Set-TextValue 'E21' '  -8.94%  '
Set-TextValue 'D22' '352.36'
Set-TextValue 'E22' '  -8.16%  '
Set-TextValue 'D23' '0.552'
Set-TextValue 'E23' '  -2.61%  '
Set-TextValue 'E24' '  -0.06%  '
Set-TextValue 'D25' '3.419.13'
Set-TextValue 'E25' '  -5.33%  '
Set-TextValue 'D26' '69.32'
Set-TextValue 'E26' '  -7.46%  '
Set-TextValue 'D28' '0.996'
Set-TextValue 'E28' '  -0.29%  '
Set-TextValue 'D29' '7.21'
Set-TextValue 'E29' '  +2.33%  '
Set-TextValue 'D30' '1.43'
Set-TextValue 'E30' '  +0.02%  '
Set-TextValue 'D31' '7.80'
Set-TextValue 'E31' '  -1.54%  '
Set-TextValue 'E32' '  -5.69%  '
Set-TextValue 'B33' 'Kaspa'
Set-TextValue 'C33' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D33' '0.150'
Set-TextValue 'E33' '  -1.84%  '
Set-TextValue 'B34' 'USDe'
Set-TextValue 'C34' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D34' '1.00'
Set-TextValue 'E34' '  +0.01%  '
Set-TextValue 'D35' '3.321.15'
Set-TextValue 'E35' '  -5.07%  '
Set-TextValue 'D36' '22.68'
Set-TextValue 'E36' '  -0.84%  '
Set-TextValue 'D37' '5.20'
Set-TextValue 'E37' '  +0.59%  '
Set-TextValue 'D38' '6.74'
Set-TextValue 'E38' '  +0.07%  '
Set-TextValue 'E39' '  -1.05%  '
Set-TextValue 'D40' '158.41'
Set-TextValue 'E40' '  -2.16%  '
Set-TextValue 'D42' '0.999'
Set-TextValue 'E42' '  -0.12%  '
Set-TextValue 'D43' '41.05'
Set-TextValue 'E43' '  -0.51%  '
Set-TextValue 'E44' '  +1.42%  '
Set-TextValue 'D45' '0.739'
Set-TextValue 'E45' '  -7.29%  '
Set-TextValue 'E46' '  +0.85%  '
Set-TextValue 'D47' '22.69'
Set-TextValue 'E47' '  -3.40%  '
Set-TextValue 'E48' '  -4.36%  '
Set-TextValue 'D49' '6.65'
Set-TextValue 'E49' '  -0.59%  '
Set-TextValue 'D50' '0.857'
Set-TextValue 'E50' '  -4.74%  '
Set-TextValue 'D51' '21.35'
Set-TextValue 'E51' '  +4.47%  '
